$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 45224 (2023-10-25) to 45233 (2023-11-03) for rows 2 through 9.
$newDate = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$newDate = $newDate.AddDays(45233)

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
